$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I51").Value = 'sd'
$ws.Range("J51").Value = 'Statement-non-opinion'
$ws.Range("I57").Value = 'sd'
$ws.Range("J57").Value = 'Statement-non-opinion'
$ws.Range("I60").Value = 'sv'
$ws.Range("J60").Value = 'Statement-opinion'
$ws.Range("I65").Value = 'sv'
$ws.Range("J65").Value = 'Statement-opinion'
$ws.Range("I79").Value = '%'
$ws.Range("J79").Value = 'Uninterpretable'
$ws.Range("I86").Value = 'aa'
$ws.Range("J86").Value = 'Agree/Accept'
$ws.Range("I103").Value = 'sd'
$ws.Range("J103").Value = 'Statement-non-opinion'
$ws.Range("I109").Value = 'sv'
$ws.Range("J109").Value = 'Statement-opinion'
$ws.Range("I163").Value = '%'
$ws.Range("J163").Value = 'Uninterpretable'
$ws.Range("I173").Value = 'sv'
$ws.Range("J173").Value = 'Statement-opinion'
$ws.Range("I183").Value = '%'
$ws.Range("J183").Value = 'Uninterpretable'
$ws.Range("I199").Value = 'aa'
$ws.Range("J199").Value = 'Agree/Accept'
$ws.Range("I202").Value = 'sv'
$ws.Range("J202").Value = 'Statement-opinion'
$ws.Range("I208").Value = 'sd'
$ws.Range("J208").Value = 'Statement-non-opinion'
$ws.Range("I228").Value = 'aa'
$ws.Range("J228").Value = 'Agree/Accept'
$ws.Range("I245").Value = 'sv'
$ws.Range("J245").Value = 'Statement-opinion'
$ws.Range("I250").Value = 'b'
$ws.Range("J250").Value = 'Acknowledge (Backchannel)'
$ws.Range("I254").Value = 'sv'
$ws.Range("J254").Value = 'Statement-opinion'
$ws.Range("I269").Value = 'sv'
$ws.Range("J269").Value = 'Statement-opinion'
$ws.Range("I271").Value = 'sv'
$ws.Range("J271").Value = 'Statement-opinion'
$ws.Range("I278").Value = 'sd'
$ws.Range("J278").Value = 'Statement-non-opinion'
$ws.Range("I288").Value = 'sv'
$ws.Range("J288").Value = 'Statement-opinion'
$ws.Range("I289").Value = 'sv'
$ws.Range("J289").Value = 'Statement-opinion'
